$wb = $excel.ActiveWorkbook

# Sheet "HESD-SYSoHEbV": remove the 59 "Vintage1891".."Vintage1949" columns
# (columns B:BH), shifting the remaining vintage columns (1950-2100) left.
$wsStart = $wb.Worksheets.Item("HESD-SYSoHEbV")
$wsStart.Range("B1:BH1").EntireColumn.Delete()

# Sheet "HESD-FoHERbA": drop the now-unneeded trailing 59 age columns
# (the tail was already saturated at 1, so truncate EX:HD).
$wsFrac = $wb.Worksheets.Item("HESD-FoHERbA")
$wsFrac.Range("EX1:HD1").EntireColumn.Delete()

$wsStart.Activate()
$wb.Worksheets.Item("HESD-SYSoHEbV").Range("E25").Select()
